$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (batsman) for ownTeam and oppTeam
$ws.Range("D1:E1").EntireColumn.Insert()

# Force the numeric-looking data cells to be stored as text (matches source
# data, all t="str" rather than numeric) without touching the header row
$ws.Range("G2:K4").NumberFormat = "@"

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Row 2 data (Abu Dhabi / October 28 2020 / Mumbai won...)
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 28 2020"
$ws.Range("C2").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Josh Philippe "
$ws.Range("G2").Value = "33"
$ws.Range("H2").Value = "24"
$ws.Range("I2").Value = "4"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "137.50"

# Row 3 data (Sharjah / October 31 2020 / Sunrisers...)
$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 31 2020"
$ws.Range("C3").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Sunrisers Hyderabad"
$ws.Range("F3").Value = "Josh Philippe "
$ws.Range("G3").Value = "32"
$ws.Range("H3").Value = "31"
$ws.Range("I3").Value = "4"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "103.22"

# Row 4 data (Abu Dhabi / November 02 2020 / Capitals...)
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " November 02 2020"
$ws.Range("C4").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Delhi Capitals"
$ws.Range("F4").Value = "Josh Philippe "
$ws.Range("G4").Value = "12"
$ws.Range("H4").Value = "17"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "70.58"
